# Updated cryptos list with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row,
# and swaps the dogwifhat / ImmutableX rows (44-45) to their new rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.404.10"
$ws.Range("E2").Value = "  -2.48%  "

# Row 3
$ws.Range("D3").Value = "2.539.16"
$ws.Range("E3").Value = "  -3.73%  "

# Row 5
$ws.Range("D5").Value = "'582.27"
$ws.Range("E5").Value = "  -1.03%  "

# Row 6
$ws.Range("D6").Value = "'169.72"
$ws.Range("E6").Value = "  -2.83%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "'0.520"
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("D9").Value = "2.539.38"
$ws.Range("E9").Value = "  -3.65%  "

# Row 10
$ws.Range("E10").Value = "  -4.17%  "

# Row 11
$ws.Range("E11").Value = "  -1.36%  "

# Row 12
$ws.Range("D12").Value = "'0.353"
$ws.Range("E12").Value = "  -1.81%  "

# Row 13
$ws.Range("E13").Value = "  -0.31%  "

# Row 14
$ws.Range("D14").Value = "2.987.41"
$ws.Range("E14").Value = "  -4.26%  "

# Row 15
$ws.Range("D15").Value = "70.237.10"
$ws.Range("E15").Value = "  -2.51%  "

# Row 16
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("E16").Value = "  -5.08%  "

# Row 17
$ws.Range("D17").Value = "'25.36"
$ws.Range("E17").Value = "  -1.53%  "

# Row 18
$ws.Range("D18").Value = "2.544.69"
$ws.Range("E18").Value = "  -4.00%  "

# Row 19
$ws.Range("D19").Value = "'7.93"
$ws.Range("E19").Value = "  -0.68%  "

# Row 20
$ws.Range("D20").Value = "'11.44"
$ws.Range("E20").Value = "  -5.55%  "

# Row 21
$ws.Range("D21").Value = "'354.04"
$ws.Range("E21").Value = "  -5.50%  "

# Row 22
$ws.Range("D22").Value = "'3.96"
$ws.Range("E22").Value = "  -2.57%  "

# Row 23
$ws.Range("D23").Value = "'1.99"
$ws.Range("E23").Value = "  -1.79%  "

# Row 24
$ws.Range("E24").Value = "  +0.09%  "

# Row 25
$ws.Range("D25").Value = "'70.03"
$ws.Range("E25").Value = "  -1.79%  "

# Row 26
$ws.Range("D26").Value = "'4.04"
$ws.Range("E26").Value = "  -3.82%  "

# Row 27
$ws.Range("D27").Value = "'9.11"
$ws.Range("E27").Value = "  -2.51%  "

# Row 28
$ws.Range("D28").Value = "2.677.58"
$ws.Range("E28").Value = "  -3.51%  "

# Row 29
$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -0.46%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0919"
$ws.Range("E30").Value = "  -3.01%  "

# Row 31
$ws.Range("D31").Value = "'7.92"
$ws.Range("E31").Value = "  -0.28%  "

# Row 32
$ws.Range("D32").Value = "'1.27"
$ws.Range("E32").Value = "  -2.53%  "

# Row 33
$ws.Range("D33").Value = "'469.70"
$ws.Range("E33").Value = "  -3.97%  "

# Row 34
$ws.Range("D34").Value = "'1.77"
$ws.Range("E34").Value = "  -1.80%  "

# Row 35
$ws.Range("E35").Value = "  -0.03%  "

# Row 36
$ws.Range("D36").Value = "'0.120"
$ws.Range("E36").Value = "  +2.80%  "

# Row 37
$ws.Range("D37").Value = "'154.97"
$ws.Range("E37").Value = "  -3.90%  "

# Row 38
$ws.Range("D38").Value = "'19.03"
$ws.Range("E38").Value = "  +0.58%  "

# Row 39
$ws.Range("D39").Value = "'18.59"
$ws.Range("E39").Value = "  -3.54%  "

# Row 40
$ws.Range("E40").Value = "  +0.01%  "

# Row 41
$ws.Range("D41").Value = "'4.83"
$ws.Range("E41").Value = "  -0.58%  "

# Row 42
$ws.Range("D42").Value = "'0.323"
$ws.Range("E42").Value = "  -1.14%  "

# Row 43
$ws.Range("E43").Value = "  -6.65%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.35"
$ws.Range("E44").Value = "  -8.92%  "

# Row 45
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "'1.17"
$ws.Range("E45").Value = "  -13.64%  "

# Row 46
$ws.Range("D46").Value = "'38.29"
$ws.Range("E46").Value = "  -2.00%  "

# Row 47
$ws.Range("D47").Value = "'144.65"
$ws.Range("E47").Value = "  -3.89%  "

# Row 48
$ws.Range("D48").Value = "'0.535"
$ws.Range("E48").Value = "  -1.06%  "

# Row 49
$ws.Range("D49").Value = "'3.53"
$ws.Range("E49").Value = "  -3.02%  "

# Row 50
$ws.Range("D50").Value = "'1.61"
$ws.Range("E50").Value = "  -3.15%  "

# Row 51
$ws.Range("D51").Value = "'0.0740"
$ws.Range("E51").Value = "  -0.55%  "
